# Applies the diff described in the commit "cambios antes de revision".
#
# Strategy: Range.InsertXML, called on the Range of an *entire existing
# paragraph*, replaces that paragraph with exactly the <w:p>...</w:p>
# fragment(s) supplied (no inherited pPr/rPr, no leftover <w:proofErr/>
# or <w:tab/> artifacts). To both edit a paragraph's own runs *and*
# splice in brand new paragraphs next to it, we pass a *sequence* of
# <w:p> elements - the original (possibly edited) paragraph plus any
# new ones - in one InsertXML call.
#
# We work from the bottom of the document upward so paragraph indices
# for content above the current edit point never shift under us.

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$wrNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph 19 ("Las tecnologías...Irrlicht Engine.") — keep its
#    text as-is, and append a new paragraph right after it:
#    "También se llevará a cabo un seguimiento de versiones mediante GIT."
# ---------------------------------------------------------------------
$p19 = $d.Paragraphs(19)
Write-Host "P19 check: [$($p19.Range.Text)]"
$r19 = $p19.Range
$xml = "<w:p $wNs><w:r><w:t xml:space=`"preserve`">Las tecnologías que vamos a aplicar para la elaboración del proyecto son principalmente C++20, como la base de todo el motor ECS, complementado con un motor gráfico para obtener una demo sencilla en dos dimensiones. Este motor será </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Irrlicht</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Engine</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t>.</w:t></w:r></w:p>" + `
       "<w:p $wNs><w:r><w:t>También se llevará a cabo un seguimiento de versiones mediante GIT.</w:t></w:r></w:p>"
$r19.InsertXML($xml)

# ---------------------------------------------------------------------
# 2) Paragraph 18 ("`tTECNOLOGÍAS PARA EL DESARROLLO") — drop the
#    leading <w:tab/> run and add <w:lastRenderedPageBreak/> before the
#    text run.
# ---------------------------------------------------------------------
$p18 = $d.Paragraphs(18)
Write-Host "P18 check: [$($p18.Range.Text)]"
$r18 = $p18.Range
$xml = "<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>TECNOLOGÍAS PARA EL DESARROLLO</w:t></w:r></w:p>"
$r18.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) Insert an extra empty paragraph right after paragraph 17 (the
#    blank line right before the TECNOLOGÍAS heading) by replacing it
#    with two empty paragraphs.
# ---------------------------------------------------------------------
$p17 = $d.Paragraphs(17)
Write-Host "P17 check: [$($p17.Range.Text)]"
$r17 = $p17.Range
$xml = "<w:p $wNs/><w:p $wNs/>"
$r17.InsertXML($xml)

# ---------------------------------------------------------------------
# 4) Paragraph 12 ("Según el autor...") — split the final sentence into
#    its own run with new wording.
# ---------------------------------------------------------------------
$p12 = $d.Paragraphs(12)
Write-Host "P12 check: [$($p12.Range.Text)]"
$r12 = $p12.Range
$xml = "<w:p $wNs><w:r><w:t>Según el autor, en este artículo no se pretende empezar y acabar un ECS completo si no explicar conceptos y obtener un resultado suficiente para su propósito</w:t></w:r><w:r><w:t>, es decir, las explicaciones que otorga son en base a un proyecto existente propio.</w:t></w:r></w:p>"
$r12.InsertXML($xml)

# ---------------------------------------------------------------------
# 5) Paragraph 8 ("El objetivo de esta pequeña guía...") — split into
#    three runs, swapping "pequeño" for "simple".
# ---------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
Write-Host "P8 check: [$($p8.Range.Text)]"
$r8 = $p8.Range
$xml = "<w:p $wNs><w:r><w:t xml:space=`"preserve`">El objetivo de esta pequeña guía es que entiendas y aprendas como funciona un </w:t></w:r><w:r><w:t>simple</w:t></w:r><w:r><w:t xml:space=`"preserve`"> ECS acoplable a diferentes tipos de juegos.</w:t></w:r></w:p>"
$r8.InsertXML($xml)

# ---------------------------------------------------------------------
# 6) Paragraph 6 (hyperlink "https://www.david-colson.com/...") — split
#    the link text into three runs while keeping the hyperlink wrapper
#    and its Hyperlink character style intact. (InsertXML silently
#    drops <w:rStyle> from supplied rPr, so the plain runs are inserted
#    first and the Hyperlink character style is re-applied afterwards
#    via Range.Style.)
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
Write-Host "P6 check: [$($p6.Range.Text)]"
$r6 = $p6.Range
$xml = "<w:p $wrNs><w:hyperlink r:id=`"rId4`" w:history=`"1`"><w:r><w:t>https://www.david-c</w:t></w:r><w:r><w:t>o</w:t></w:r><w:r><w:t>lson.com/2020/02/09/making-a-simple-ecs.html</w:t></w:r></w:hyperlink></w:p>"
$r6.InsertXML($xml)
$p6b = $d.Paragraphs(6)
$r6b = $p6b.Range
$r6b.MoveEnd(1, -1)
$r6b.Style = "Hyperlink"

# ---------------------------------------------------------------------
# 7) Paragraph 3 ("`tANTECEDENTES") — drop the leading <w:tab/> run.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
Write-Host "P3 check: [$($p3.Range.Text)]"
$r3 = $p3.Range
$xml = "<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ANTECEDENTES</w:t></w:r></w:p>"
$r3.InsertXML($xml)

# ---------------------------------------------------------------------
# 8) Paragraph 2 ("Intro ") — replace with three new paragraphs: two
#    plain intro paragraphs plus the (bold-marked) replacement for
#    "Intro ".
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
Write-Host "P2 check: [$($p2.Range.Text)]"
$r2 = $p2.Range
$xml = "<w:p $wNs><w:r><w:t>En el apartado actual, es difícil catalogar como igual a otros proyectos existentes acerca de como hacer un videojuego o pequeños post o webs que resuelvan problemas cercanos, pero no el mismo.</w:t></w:r></w:p>" + `
       "<w:p $wNs><w:r><w:t>En el siguiente apartado he comentado tres " + [char]0x201C + "proyectos" + [char]0x201D + " que tratan sobre la enseñanza de como hacer un ECS, de una forma incompleta, o basados en proyectos ya existentes.</w:t></w:r></w:p>" + `
       "<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>Mi búsqueda de arte similar no ha obtenido resultados, es por esto por lo que no tengo una referencia de proyecto similar.</w:t></w:r></w:p>"
$r2.InsertXML($xml)

Write-Host "DONE"
